# Update LR-pairs sheet with newly computed TPM-based values.
# Rows 5:7 (previously duplicate "MuSCs" sending-cluster rows) are removed,
# and rows 2:4 are refreshed with the new MuSCs-based TPM figures
# (replacing the old "ECs" sending-cluster rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 5-7 first (shifts nothing else, data kept in 2-4).
$ws.Rows("5:7").Delete()

# Row 2: MuSCs -> Il10 -> Il10ra -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Il10"
$ws.Range("C2").Value = "Il10ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2168213333333333
$ws.Range("H2").Value = 0.650464
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07334
$ws.Range("N2").Value = 0.22002
$ws.Range("O2").Value = 0.7205147937713883
$ws.Range("P2").Value = 0.7205147937713883
$ws.Range("Q2").Value = 0.01590167658666667
$ws.Range("R2").Value = 0.14311508928
$ws.Range("S2").Value = 0.7205147937713883
$ws.Range("T2").Value = 0.7205147937713883

# Row 3: MuSCs -> Il10 -> Il10ra -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Il10"
$ws.Range("C3").Value = "Il10ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2168213333333333
$ws.Range("H3").Value = 0.650464
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01293133333333333
$ws.Range("N3").Value = 0.038794
$ws.Range("O3").Value = 0.1270414094608092
$ws.Range("P3").Value = 0.1270414094608092
$ws.Range("Q3").Value = 0.002803788935111112
$ws.Range("R3").Value = 0.025234100416
$ws.Range("S3").Value = 0.1270414094608092
$ws.Range("T3").Value = 0.1270414094608092

# Row 4: MuSCs -> Il10 -> Il10ra -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Il10"
$ws.Range("C4").Value = "Il10ra"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2168213333333333
$ws.Range("H4").Value = 0.650464
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.015517
$ws.Range("N4").Value = 0.046551
$ws.Range("O4").Value = 0.1524437967678025
$ws.Range("P4").Value = 0.1524437967678025
$ws.Range("Q4").Value = 0.003364416629333334
$ws.Range("R4").Value = 0.030279749664
$ws.Range("S4").Value = 0.1524437967678025
$ws.Range("T4").Value = 0.1524437967678025

$wb.Save()
